$d = $word.ActiveDocument

# Locate the sentence that needs the word "still" inserted into it, right
# after "...does not wish to be married may" and before " be made Partners."
$rng = $d.Content
$found = $rng.Find.Execute(
    "any pairing who does not wish to be married may",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "", 0
)

if ($found) {
    # Collapse to the insertion point right after "...may" and insert the
    # missing word. Track the insertion as a revision and then accept just
    # that revision -- this preserves the run split (prefix run / " still"
    # run / suffix run) that a genuine in-place edit produces, instead of
    # silently re-merging the text back into a single run.
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true

    $rng.Collapse(0)
    $rng.InsertAfter(" still")

    for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
        $d.Revisions($i).Accept()
    }

    $d.TrackRevisions = $wasTracking
}
